# "Minha parte do relatório"
# Applies:
#   - Footer/date placeholder cached text update (all slide layouts + slide master)
#   - Slide 9 layout rework: move several shapes, delete two images/connectors,
#     resize + reword one callout box.

$p = $ppt.ActivePresentation
$EMU = 12700.0  # points -> EMU

# ---------------------------------------------------------------------------
# 1) Date placeholder ("01/09/2017" -> "02/12/2017") cached field text, on the
#    slide master and on every slide layout.
# ---------------------------------------------------------------------------
function Update-DatePlaceholder($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shp = $shapes.Item($i)
        if ($shp.Name -like "Date Placeholder*") {
            $shp.TextFrame.TextRange.Text = "02/12/2017"
        }
    }
}

$master = $p.SlideMaster
Update-DatePlaceholder $master.Shapes

for ($li = 1; $li -le $master.CustomLayouts.Count; $li++) {
    $layout = $master.CustomLayouts.Item($li)
    Update-DatePlaceholder $layout.Shapes
}

# ---------------------------------------------------------------------------
# 2) Slide 9 edits
# ---------------------------------------------------------------------------
$s = $p.Slides.Item(9)

# Retângulo 68: reposition
$shp = $s.Shapes.Item("Retângulo 68")
$shp.Left = 3736266 / $EMU
$shp.Top = 4336017 / $EMU

# Remove "Imagem 22" (duplicate vector-graphic icon)
$s.Shapes.Item("Imagem 22").Delete()

# Imagem 24: reposition
$shp = $s.Shapes.Item("Imagem 24")
$shp.Left = 4135009 / $EMU
$shp.Top = 4557742 / $EMU

# Remove the two straight-arrow connectors that pointed at the deleted icon
$s.Shapes.Item("Conector de Seta Reta 52").Delete()
$s.Shapes.Item("Conector de Seta Reta 53").Delete()

# Conector de Seta Reta 54: reposition
$shp = $s.Shapes.Item("Conector de Seta Reta 54")
$shp.Left = 4580988 / $EMU
$shp.Top = 3825475 / $EMU

# Conector de Seta Reta 55: reposition
$shp = $s.Shapes.Item("Conector de Seta Reta 55")
$shp.Left = 4750831 / $EMU
$shp.Top = 3850937 / $EMU

# CaixaDeTexto 71 ("Aplicativo/Web"): reposition
$shp = $s.Shapes.Item("CaixaDeTexto 71")
$shp.Left = 3196661 / $EMU
$shp.Top = 5600362 / $EMU

# CaixaDeTexto 79: grow + reword
$shp = $s.Shapes.Item("CaixaDeTexto 79")
$shp.Height = 1015663 / $EMU
$shp.TextFrame.TextRange.Text = "Processa a informação e classifica o local, tipo e disponibilidade da vaga na região e dados do Usuário"

# CaixaDeTexto 80: reposition
$shp = $s.Shapes.Item("CaixaDeTexto 80")
$shp.Left = 4335599 / $EMU
$shp.Top = 1769296 / $EMU

# CaixaDeTexto 82: reposition
$shp = $s.Shapes.Item("CaixaDeTexto 82")
$shp.Left = 2535263 / $EMU
$shp.Top = 5260802 / $EMU
